# Apply the 2020-10-25 data refresh for "Fonds de solidarite volet 2" (regional x categorie
# juridique). The source extract grew from 105 to 107 data rows; most existing nombre_aides /
# montant_total figures were revised upward, and the per-region categorie-juridique breakdown
# shifted for Occitanie / Provence-Alpes-Cote d'Azur (extra categories appear in the new extract).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the column is one of the numeric-looking
# columns (nombre_aides, montant_total, reg, code_categorie_juridique) that must stay stored
# as text, matching the source file's inlineStr cells (Excel would otherwise auto-convert a
# bare numeric-looking string into a real Number cell).
$updates = @(
    ,@("C2", "1694", $true)
    ,@("D2", "4226936.79", $true)
    ,@("C4", "1259", $true)
    ,@("D4", "6405724.18", $true)
    ,@("C6", "867", $true)
    ,@("D6", "3896453.05", $true)
    ,@("C8", "48", $true)
    ,@("D8", "221070.31", $true)
    ,@("C9", "257", $true)
    ,@("D9", "841266.05", $true)
    ,@("C10", "5", $true)
    ,@("D10", "19000.00", $true)
    ,@("C11", "448", $true)
    ,@("D11", "2570799.28", $true)
    ,@("C12", "215", $true)
    ,@("D12", "1325466.34", $true)
    ,@("C15", "283", $true)
    ,@("D15", "982376.49", $true)
    ,@("C17", "640", $true)
    ,@("D17", "5158888.72", $true)
    ,@("C18", "182", $true)
    ,@("D18", "1209225.93", $true)
    ,@("C20", "11", $true)
    ,@("D20", "44309.04", $true)
    ,@("C21", "227", $true)
    ,@("D21", "698413.14", $true)
    ,@("C23", "430", $true)
    ,@("D23", "2498187.61", $true)
    ,@("C24", "195", $true)
    ,@("D24", "1038686.27", $true)
    ,@("C26", "13", $true)
    ,@("D26", "58250.00", $true)
    ,@("C31", "414", $true)
    ,@("D31", "1208791.11", $true)
    ,@("C32", "10", $true)
    ,@("D32", "50000.00", $true)
    ,@("C33", "782", $true)
    ,@("D33", "5133127.55", $true)
    ,@("C34", "4", $true)
    ,@("D34", "15500.00", $true)
    ,@("C35", "510", $true)
    ,@("D35", "2711240.28", $true)
    ,@("C37", "21", $true)
    ,@("D37", "75432.00", $true)
    ,@("C38", "567", $true)
    ,@("D38", "1496968.32", $true)
    ,@("C39", "261", $true)
    ,@("D39", "912292.04", $true)
    ,@("C40", "263", $true)
    ,@("D40", "883238.72", $true)
    ,@("C43", "397", $true)
    ,@("D43", "1397718.40", $true)
    ,@("C44", "177", $true)
    ,@("D44", "1016977.39", $true)
    ,@("C45", "237", $true)
    ,@("D45", "1100076.19", $true)
    ,@("C47", "14", $true)
    ,@("D47", "82085.23", $true)
    ,@("C48", "701", $true)
    ,@("D48", "2038180.93", $true)
    ,@("C50", "940", $true)
    ,@("D50", "5863690.56", $true)
    ,@("C51", "671", $true)
    ,@("D51", "3555468.96", $true)
    ,@("C54", "8767", $true)
    ,@("D54", "22752081.68", $true)
    ,@("C57", "42", $true)
    ,@("D57", "307600.00", $true)
    ,@("C58", "6049", $true)
    ,@("D58", "28120428.74", $true)
    ,@("C59", "20", $true)
    ,@("D59", "230000.00", $true)
    ,@("C60", "5935", $true)
    ,@("D60", "23625494.17", $true)
    ,@("C61", "64", $true)
    ,@("D61", "174170.00", $true)
    ,@("C62", "123", $true)
    ,@("D62", "525315.40", $true)
    ,@("C67", "175", $true)
    ,@("D67", "444150.00", $true)
    ,@("C68", "248", $true)
    ,@("D68", "737482.59", $true)
    ,@("C69", "183", $true)
    ,@("D69", "468744.35", $true)
    ,@("C70", "9", $true)
    ,@("D70", "28421.00", $true)
    ,@("C71", "117", $true)
    ,@("D71", "487784.00", $true)
    ,@("C74", "278", $true)
    ,@("D74", "948768.30", $true)
    ,@("C75", "482", $true)
    ,@("D75", "2941810.53", $true)
    ,@("C76", "290", $true)
    ,@("D76", "2020164.29", $true)
    ,@("C79", "437", $true)
    ,@("D79", "1294368.80", $true)
    ,@("C80", "7", $true)
    ,@("D80", "42254.00", $true)
    ,@("C81", "1173", $true)
    ,@("D81", "7608638.94", $true)
    ,@("C82", "625", $true)
    ,@("D82", "3517220.36", $true)
    ,@("C84", "42", $true)
    ,@("D84", "189227.77", $true)
    ,@("C85", "779", $true)
    ,@("D85", "2086083.64", $true)
    ,@("C88", "1125", $true)
    ,@("D88", "5147824.50", $true)
    ,@("C89", "3", $true)
    ,@("D89", "57000.00", $true)
    ,@("G89", "55", $true)
    ,@("H89", "Société anonyme à conseil d'administration", $false)
    ,@("C90", "814", $true)
    ,@("D90", "3493559.88", $true)
    ,@("G90", "57", $true)
    ,@("H90", "Société par actions simplifiée", $false)
    ,@("C91", "33", $true)
    ,@("D91", "91206.00", $true)
    ,@("G91", "65", $true)
    ,@("H91", "Société civile", $false)
    ,@("C92", "36", $true)
    ,@("D92", "133571.23", $true)
    ,@("E92", "76", $true)
    ,@("F92", "Occitanie", $false)
    ,@("G92", "92", $true)
    ,@("H92", "Association loi 1901 ou assimilé", $false)
    ,@("C93", "248", $true)
    ,@("D93", "591350.00", $true)
    ,@("G93", "10", $true)
    ,@("H93", "Entrepreneur individuel", $false)
    ,@("C94", "4", $true)
    ,@("D94", "9500.00", $true)
    ,@("G94", "52", $true)
    ,@("H94", "Société en nom collectif", $false)
    ,@("C95", "607", $true)
    ,@("D95", "3509341.95", $true)
    ,@("G95", "54", $true)
    ,@("H95", "Société à responsabilité limitée (SARL)", $false)
    ,@("C96", "3", $true)
    ,@("D96", "60000.00", $true)
    ,@("G96", "55", $true)
    ,@("H96", "Société anonyme à conseil d'administration", $false)
    ,@("C97", "217", $true)
    ,@("D97", "911096.11", $true)
    ,@("G97", "57", $true)
    ,@("H97", "Société par actions simplifiée", $false)
    ,@("C98", "18", $true)
    ,@("D98", "60500.00", $true)
    ,@("E98", "52", $true)
    ,@("F98", "Pays de la Loire", $false)
    ,@("G98", "65", $true)
    ,@("H98", "Société civile", $false)
    ,@("C99", "9", $true)
    ,@("D99", "34670.00", $true)
    ,@("E99", "52", $true)
    ,@("F99", "Pays de la Loire", $false)
    ,@("G99", "92", $true)
    ,@("H99", "Association loi 1901 ou assimilé", $false)
    ,@("C100", "1319", $true)
    ,@("D100", "3363754.28", $true)
    ,@("G100", "10", $true)
    ,@("H100", "Entrepreneur individuel", $false)
    ,@("C101", "3", $true)
    ,@("D101", "7571.16", $true)
    ,@("G101", "22", $true)
    ,@("H101", "Société créée de fait", $false)
    ,@("C102", "12", $true)
    ,@("D102", "47560.00", $true)
    ,@("G102", "52", $true)
    ,@("H102", "Société en nom collectif", $false)
    ,@("C103", "1527", $true)
    ,@("D103", "7250036.93", $true)
    ,@("G103", "54", $true)
    ,@("H103", "Société à responsabilité limitée (SARL)", $false)
    ,@("C104", "3", $true)
    ,@("D104", "7500.00", $true)
    ,@("G104", "55", $true)
    ,@("H104", "Société anonyme à conseil d'administration", $false)
    ,@("C105", "1469", $true)
    ,@("D105", "6230699.20", $true)
    ,@("G105", "57", $true)
    ,@("H105", "Société par actions simplifiée", $false)
    ,@("A106", "Fonds de solidarité", $false)
    ,@("B106", "VOLET2", $false)
    ,@("C106", "19", $true)
    ,@("D106", "55345.28", $true)
    ,@("E106", "93", $true)
    ,@("F106", "Provence-Alpes-Côte d'Azur", $false)
    ,@("G106", "65", $true)
    ,@("H106", "Société civile", $false)
    ,@("A107", "Fonds de solidarité", $false)
    ,@("B107", "VOLET2", $false)
    ,@("C107", "80", $true)
    ,@("D107", "360788.23", $true)
    ,@("E107", "93", $true)
    ,@("F107", "Provence-Alpes-Côte d'Azur", $false)
    ,@("G107", "92", $true)
    ,@("H107", "Association loi 1901 ou assimilé", $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $value = $u[1]
    $isNumericText = $u[2]
    $cell = $ws.Range($ref)
    if ($isNumericText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

